$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$sa = $shp.SmartArt
$nodes = $sa.AllNodes
$n = $nodes.Item(3)
$tf = $n.TextFrame
$tr = $tf.TextRange
$tr.Text = "O token recebido pode ficar guardado em memória, até que o browser/tab seja fechado ou seja feito um logout (limpar o token). Alternativamente, o token pode ser guardado no session storage (sessionStorage.setItem), para não se perder, em caso de refresh."
